$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136, shifting rows 136-174 down to 137-175
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with data.
# Static columns copied from the template row (same across this block of rows).
$ws.Cells.Item(136, 1).Value = 7
$ws.Cells.Item(136, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(136, 3).Value = "Ñuble"
$ws.Cells.Item(136, 4).Value = 44463
$ws.Cells.Item(136, 5).Value = 16
$ws.Cells.Item(136, 6).Value = 100114013
$ws.Cells.Item(136, 7).Value = "Zanahoria"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 160
$ws.Cells.Item(136, 11).Value = 6500
$ws.Cells.Item(136, 12).Value = 7000
$ws.Cells.Item(136, 13).Value = 6750
$ws.Cells.Item(136, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(136, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(136, 16).Value = 338
$ws.Cells.Item(136, 17).Value = 20
$ws.Cells.Item(136, 18).Value = "Hortaliza"
